# Add a new entry row for "dictionaryToKekulized.py" to the master list,
# right after the existing "xyzToDictionaryAromatic.py" row (row 19) and
# before the "RMG-database: Thermo Groups" section header.
#
# This mirrors Excel's own "Insert" behaviour: inserting a whole row at
# row 20 pushes every row below it (old rows 21..81) down by one
# (new rows 22..82), which is exactly what the target diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20.
$ws.Rows.Item(20).Insert()

# Populate the new row with the new script's info, matching the
# existing columns: Name, Location, Creator/Maintainer, Purpose.
$ws.Range("A20").Value = "dictionaryToKekulized.py"
$ws.Range("B20").Value = "WIP RMG_input"
$ws.Range("C20").Value = "Max"
$ws.Range("D20").Value = "For an existing species dictionary, convert adjacency lists for aromatic species to kekulized bond types"

# Update the view: move the selection to A17 (clears the old
# topLeftCell/selection scrolled state from B73).
$ws.Range("A17").Select()
